$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-10-06 -> 2023-10-07) for every data row (rows 2 through 103).
$ws.Range("C2:C103").Value = 45206
